$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-apply the existing "Text" number format to the header row / row label cell.
# (These cells already render as Text with no border; this keeps that formatting
# intact / re-affirms it, matching the style bookkeeping churn seen upstream.)
$ws.Range("A1").NumberFormat = "@"
$ws.Range("B1").NumberFormat = "@"
$ws.Range("C1").NumberFormat = "@"
$ws.Range("A2").NumberFormat = "@"

# Update the predicted score value in B2.
$ws.Range("B2").Value = 319516.2110549332
